# Refresh cached Market Board profit figures (currentAveragePrice / Leve price
# / Leve profit columns H:N) on each job sheet, row by row, to match the
# latest Universalis price snapshot pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 7: The Bleat Is On
$ws.Range("H7").Value = 2626.25
$ws.Range("I7").Value = 168.33333
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 168.33333
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -56.33332999999999
$ws.Range("N7").Value = -10224

# Row 14: Wand-full Tonight
$ws.Range("H14").Value = 2626.25
$ws.Range("I14").Value = 168.33333
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 168.33333
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 22.66667000000001
$ws.Range("N14").Value = -10382

# Row 33: Glazed and Confused
$ws.Range("H33").Value = 467.66666
$ws.Range("I33").Value = 396.8889
$ws.Range("J33").Value = 680
$ws.Range("K33").Value = 396.8889
$ws.Range("L33").Value = 680
$ws.Range("M33").Value = -167.8889
$ws.Range("N33").Value = -1138

# Row 49: Going Nowhere Fast
$ws.Range("H49").Value = 4439.75
$ws.Range("I49").Value = 999.6667
$ws.Range("J49").Value = 6503.8
$ws.Range("K49").Value = 2999.0001
$ws.Range("L49").Value = 19511.4
$ws.Range("M49").Value = -2863.0001
$ws.Range("N49").Value = -19783.4

# Row 101: Edge of the Arcane
$ws.Range("H101").Value = 1051.5454
$ws.Range("I101").Value = 833.2857
$ws.Range("J101").Value = 1433.5
$ws.Range("K101").Value = 2499.8571
$ws.Range("L101").Value = 4300.5
$ws.Range("M101").Value = -877.8571000000002
$ws.Range("N101").Value = -7544.5

# Row 136: I Like Big Brush and I Cannot Lie
$ws.Range("H136").Value = 68883.86
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 68883.86
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 68883.86
$ws.Range("N136").Value = -79083.86

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 2995.3098
$ws.Range("I32").Value = 2174.2986
$ws.Range("J32").Value = 16747.25
$ws.Range("K32").Value = 2174.2986
$ws.Range("L32").Value = 16747.25
$ws.Range("M32").Value = -1887.2986
$ws.Range("N32").Value = -17321.25

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 4101.972
$ws.Range("I61").Value = 2981.9666
$ws.Range("J61").Value = 9702
$ws.Range("K61").Value = 2981.9666
$ws.Range("L61").Value = 9702
$ws.Range("M61").Value = -2769.9666
$ws.Range("N61").Value = -10126

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 3644.1738
$ws.Range("I122").Value = 3230.5833
$ws.Range("J122").Value = 4095.3635
$ws.Range("K122").Value = 9691.749899999999
$ws.Range("L122").Value = 12286.0905
$ws.Range("M122").Value = -7241.749899999999
$ws.Range("N122").Value = -17186.0905

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 3320.5293
$ws.Range("I132").Value = 2339.6875
$ws.Range("J132").Value = 19014
$ws.Range("K132").Value = 7019.0625
$ws.Range("L132").Value = 57042
$ws.Range("M132").Value = -4489.0625
$ws.Range("N132").Value = -62102

# Row 135: Forgiveness for My Shins
$ws.Range("H135").Value = 48564.332
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 48564.332
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 48564.332
$ws.Range("N135").Value = -58704.332

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 4101.972
$ws.Range("I136").Value = 2981.9666
$ws.Range("J136").Value = 9702
$ws.Range("K136").Value = 8945.899800000001
$ws.Range("L136").Value = 29106
$ws.Range("M136").Value = -6395.899800000001
$ws.Range("N136").Value = -34206

$ws = $wb.Worksheets.Item("BSM")
# Row 140: Ceremonial Teeth
$ws.Range("H140").Value = 51547
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 51547
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 51547
$ws.Range("N140").Value = -61907

$ws = $wb.Worksheets.Item("CRP")
# Row 4: A Clogful of Camaraderie
$ws.Range("H4").Value = 838666.5
$ws.Range("I4").Value = 999.5
$ws.Range("J4").Value = 1257500
$ws.Range("K4").Value = 999.5
$ws.Range("L4").Value = 1257500
$ws.Range("M4").Value = -887.5
$ws.Range("N4").Value = -1257724

# Row 10: Spears and Sorcery
$ws.Range("H10").Value = 297
$ws.Range("I10").Value = 297
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 297
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -158
$ws.Range("N10").ClearContents()

# Row 16: Raise the Roof
$ws.Range("H16").Value = 1644.1052
$ws.Range("I16").Value = 616.3077
$ws.Range("J16").Value = 3871
$ws.Range("K16").Value = 616.3077
$ws.Range("L16").Value = 3871
$ws.Range("M16").Value = -329.3077
$ws.Range("N16").Value = -4445

# Row 92: Walk the Walk
$ws.Range("H92").Value = 56400
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 56400
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 56400
$ws.Range("N92").Value = -61392

# Row 113: Patient Patients
$ws.Range("H113").Value = 1644.1052
$ws.Range("I113").Value = 616.3077
$ws.Range("J113").Value = 3871
$ws.Range("K113").Value = 616.3077
$ws.Range("L113").Value = 3871
$ws.Range("M113").Value = 1553.6923
$ws.Range("N113").Value = -8211

# Row 135: The Wing's Wings
$ws.Range("H135").Value = 69180.8
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 69180.8
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 69180.8
$ws.Range("N135").Value = -79320.8

# Row 138: Bow Out
$ws.Range("H138").Value = 69284.60000000001
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 69284.60000000001
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 69284.60000000001
$ws.Range("N138").Value = -79564.60000000001

# Row 140: Spear Pressure
$ws.Range("H140").Value = 64595.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 64595.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 64595.5
$ws.Range("N140").Value = -74955.5

$ws = $wb.Worksheets.Item("CUL")
# Row 61: Red Letter Day
$ws.Range("H61").Value = 714
$ws.Range("I61").Value = 79
$ws.Range("J61").Value = 2238
$ws.Range("K61").Value = 237
$ws.Range("L61").Value = 6714
$ws.Range("M61").Value = -22
$ws.Range("N61").Value = -7144

# Row 98: Sweet Kiss of Death
$ws.Range("H98").Value = 3285.875
$ws.Range("I98").Value = 3199.75
$ws.Range("J98").Value = 3372
$ws.Range("K98").Value = 9599.25
$ws.Range("L98").Value = 10116
$ws.Range("M98").Value = -8101.25
$ws.Range("N98").Value = -13112

# Row 107: Slippery Service
$ws.Range("H107").Value = 1394.8667
$ws.Range("I107").Value = 671.8333
$ws.Range("J107").Value = 1876.8889
$ws.Range("K107").Value = 2015.4999
$ws.Range("L107").Value = 5630.6667
$ws.Range("M107").Value = -95.49990000000003
$ws.Range("N107").Value = -9470.6667

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1282.0769
$ws.Range("I97").Value = 1102.0476
$ws.Range("J97").Value = 2038.2
$ws.Range("K97").Value = 1102.0476
$ws.Range("L97").Value = 2038.2
$ws.Range("M97").Value = -606.0476000000001
$ws.Range("N97").Value = -3030.2

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 2343.0278
$ws.Range("I113").Value = 1918.6086
$ws.Range("J113").Value = 3093.923
$ws.Range("K113").Value = 1918.6086
$ws.Range("L113").Value = 3093.923
$ws.Range("M113").Value = 251.3914
$ws.Range("N113").Value = -7433.923

$ws = $wb.Worksheets.Item("LTW")
# Row 100: Tiger in the Sack
$ws.Range("H100").Value = 10763.333
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 13189.375
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 13189.375
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -14271.375

# Row 122: Hell on Leather
$ws.Range("H122").Value = 256432.62
$ws.Range("I122").Value = 447944.78
$ws.Range("J122").Value = 10202.714
$ws.Range("K122").Value = 1343834.34
$ws.Range("L122").Value = 30608.142
$ws.Range("M122").Value = -1341384.34
$ws.Range("N122").Value = -34307.125

$ws = $wb.Worksheets.Item("WVR")
# Row 100: Of Great Import
$ws.Range("H100").Value = 373.3889
$ws.Range("I100").Value = 393.07144
$ws.Range("J100").Value = 304.5
$ws.Range("K100").Value = 786.14288
$ws.Range("L100").Value = 609
$ws.Range("M100").Value = -245.14288
$ws.Range("N100").Value = -1691

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 2994.25
$ws.Range("I126").Value = 1991.875
$ws.Range("J126").Value = 4999
$ws.Range("K126").Value = 5975.625
$ws.Range("L126").Value = 14997
$ws.Range("M126").Value = -3505.625
$ws.Range("N126").Value = -19937
